$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.376.00"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "2.596.56"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'587.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.86%  "
$ws.Range("D6").Value = "'149.32"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.583"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("D10").Value = "'5.81"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").Value = "'0.385"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").Value = "'27.53"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "3.062.75"
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").Value = "63.202.18"
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("D17").Value = "2.608.71"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "'12.01"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").Value = "'4.65"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "'343.56"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("D21").Value = "'6.81"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'66.46"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").Value = "'1.72"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("D25").Value = "'9.17"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("D26").Value = "'1.64"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.45%  "
$ws.Range("D27").Value = "'565.11"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.08%  "
$ws.Range("D28").Value = "'8.16"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("E29").Value = "  -3.08%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "'2.01"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.00%  "
$ws.Range("D32").Value = "0.0₃0840"
$ws.Range("E32").Value = "  -2.98%  "
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "'5.28"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").Value = "'165.68"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").Value = "'0.411"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'19.29"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("E39").Value = "  -6.34%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "'165.91"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D43").Value = "'22.63"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.13%  "
$ws.Range("D44").Value = "'0.0579"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").Value = "'2.10"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.20%  "
$ws.Range("D46").Value = "'0.629"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D49").Value = "'18.98"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.31%  "
$ws.Range("D50").Value = "0.0₆0227"
$ws.Range("E50").Value = "  +13.85%  "
$ws.Range("D51").Value = "'0.178"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.68%  "
